$d = $word.ActiveDocument

# --- 1) Mint the "Tablaconcuadrcula" (Table Grid) table style in styles.xml ---
# Use a disposable scratch table at the very start of the document purely to
# attach + configure the style, then delete the scratch table. The style
# definition remains behind in styles.xml.
$scratchRange = $d.Range(0, 0)
$scratchTable = $d.Tables.Add($scratchRange, 1, 1)
$scratchTable.Style = "Tablaconcuadrcula"
$tableStyle = $scratchTable.Style
$tableStyle.NameLocal = "Table Grid"
$tableStyle.BaseStyle = "Tablanormal"
$tableStyle.Priority = 39
$styleParaFormat = $tableStyle.ParagraphFormat
$styleParaFormat.SpaceAfter = 0
$styleParaFormat.LineSpacingRule = 0
$scratchTable.Delete()

# --- 2) Insert the page-break paragraph right before the trailing empty paragraph ---
$lastPar = $d.Paragraphs.Last
$insertPos = $lastPar.Range.Start
$breakRange = $d.Range($insertPos, $insertPos)
$breakRange.InsertBreak(7)
$pageBreakParaIndex = $d.Paragraphs.Count - 1
$pageBreakPara = $d.Paragraphs($pageBreakParaIndex)
$pageBreakPara.Range.Font.Underline = 1

# --- 3) Insert the SWOT table right before the (now further shifted) trailing paragraph ---
$lastPar2 = $d.Paragraphs.Last
$tblPos = $lastPar2.Range.Start
$tableRange = $d.Range($tblPos, $tblPos)
$tableXml = @'
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:tblPr><w:tblStyle w:val="Tablaconcuadrcula"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="4414"/><w:gridCol w:w="4414"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="4414" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>Strength</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4414" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Weaknesses</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="4414" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Advanced Functionality: The ability to detect geographic relations within the text and generate a spatial-temporal graph is an advanced functionality that can significantly improve the system's capabilities.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>User Experience: It can enhance user experience by providing more context-aware and visually appealing responses.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4414" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Technical Complexity: Developing a module that accurately detects geographic relations and generates a spatial-temporal graph can be technically challenging.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Data Availability: The accuracy and effectiveness of the functionality heavily depend on the availability and quality of the geographic data.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="4414" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Opportunities</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4414" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Threats</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="4414" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Innovation: This could lead to the development of new features or products based on advanced geographic data analysis.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Market Differentiation: By prioritizing geographic relations detection, the client could differentiate themselves in the market.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4414" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Competition: There are already established players in the market offering similar functionalities.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t>User Expectations: High user expectations for accurate geographic relations detection could lead to dissatisfaction if not met.</w:t></w:r></w:p></w:tc></w:tr></w:tbl>
'@
$tableRange.InsertXML($tableXml)

Write-Host "Edit complete."
